$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I0 and IF), matching the existing
# bold/centered/bordered header style used by the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for I2:J65
$iValues = @(4,6,9,7,5,10,8,6,7,9,7,7,10,7,9,8,6,5,9,7,6,11,6,1,8,8,5,5,6,7,9,8,8,6,8,7,6,6,5,7,8,8,7,7,8,8,9,9,9,5,8,7,9,10,7,6,8,8,4,4,5,5,6,5)
$jValues = @(5,7,9,8,6,10,8,7,8,9,7,7,10,8,9,8,7,6,9,8,7,11,6,2,8,8,6,5,6,7,9,8,8,6,8,7,6,6,6,7,9,8,7,7,8,8,9,9,9,5,8,7,9,10,7,7,9,8,4,4,5,5,6,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
